$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in title cell A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 17:03"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1240654
$ws.Cells.Item(4, 3).Value = 3021
$ws.Cells.Item(4, 4).Value = 201858
$ws.Cells.Item(4, 5).Value = 966342
$ws.Cells.Item(4, 6).Value = 16179
$ws.Cells.Item(4, 7).Value = 183
$ws.Cells.Item(4, 8).Value = 72454

# Row 18: India
$ws.Cells.Item(18, 1).Value = "India"
$ws.Cells.Item(18, 2).Value = 49852
$ws.Cells.Item(18, 3).Value = 452
$ws.Cells.Item(18, 4).Value = 14367
$ws.Cells.Item(18, 5).Value = 33781
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 11
$ws.Cells.Item(18, 8).Value = 1704

# Row 45: Republica Dominicana
$ws.Cells.Item(45, 1).Value = "Republica Dominicana"
$ws.Cells.Item(45, 2).Value = 8807
$ws.Cells.Item(45, 3).Value = 327
$ws.Cells.Item(45, 4).Value = 1905
$ws.Cells.Item(45, 5).Value = 6540
$ws.Cells.Item(45, 6).Value = 144
$ws.Cells.Item(45, 7).Value = 8
$ws.Cells.Item(45, 8).Value = 362

# Row 46: Colombia
$ws.Cells.Item(46, 1).Value = "Colombia"
$ws.Cells.Item(46, 2).Value = 8613
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 2013
$ws.Cells.Item(46, 5).Value = 6222
$ws.Cells.Item(46, 6).Value = 128
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 378

# Row 80: Bulgaria
$ws.Cells.Item(80, 1).Value = "Bulgaria"
$ws.Cells.Item(80, 2).Value = 1778
$ws.Cells.Item(80, 3).Value = 74
$ws.Cells.Item(80, 4).Value = 360
$ws.Cells.Item(80, 5).Value = 1334
$ws.Cells.Item(80, 6).Value = 38
$ws.Cells.Item(80, 7).Value = 4
$ws.Cells.Item(80, 8).Value = 84

# Row 112: Maldivas
$ws.Cells.Item(112, 1).Value = "Maldivas"
$ws.Cells.Item(112, 2).Value = 618
$ws.Cells.Item(112, 3).Value = 45
$ws.Cells.Item(112, 4).Value = 20
$ws.Cells.Item(112, 5).Value = 596
$ws.Cells.Item(112, 6).Value = 2
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 2

# Row 113: Mali
$ws.Cells.Item(113, 1).Value = "Mali"
$ws.Cells.Item(113, 2).Value = 612
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 228
$ws.Cells.Item(113, 5).Value = 352
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 32

# Row 114: Georgia
$ws.Cells.Item(114, 1).Value = "Georgia"
$ws.Cells.Item(114, 2).Value = 610
$ws.Cells.Item(114, 3).Value = 6
$ws.Cells.Item(114, 4).Value = 269
$ws.Cells.Item(114, 5).Value = 332
$ws.Cells.Item(114, 6).Value = 6
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 9

# Row 115: San Marino
$ws.Cells.Item(115, 1).Value = "San Marino"
$ws.Cells.Item(115, 2).Value = 608
$ws.Cells.Item(115, 3).Value = 19
$ws.Cells.Item(115, 4).Value = 97
$ws.Cells.Item(115, 5).Value = 470
$ws.Cells.Item(115, 6).Value = 4
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 41

# Row 116: Kenia
$ws.Cells.Item(116, 1).Value = "Kenia"
$ws.Cells.Item(116, 2).Value = 582
$ws.Cells.Item(116, 3).Value = 47
$ws.Cells.Item(116, 4).Value = 190
$ws.Cells.Item(116, 5).Value = 366
$ws.Cells.Item(116, 6).Value = 2
$ws.Cells.Item(116, 7).Value = 2
$ws.Cells.Item(116, 8).Value = 26

# Row 183: Yemen
$ws.Cells.Item(183, 1).Value = "Yemen"
$ws.Cells.Item(183, 2).Value = 25
$ws.Cells.Item(183, 3).Value = 3
$ws.Cells.Item(183, 4).Value = 1
$ws.Cells.Item(183, 5).Value = 19
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 1
$ws.Cells.Item(183, 8).Value = 5

# Row 184: Antigua y Barbuda
$ws.Cells.Item(184, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(184, 2).Value = 25
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 16
$ws.Cells.Item(184, 5).Value = 6
$ws.Cells.Item(184, 6).Value = 1
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 3

# Row 185: Timor Oriental
$ws.Cells.Item(185, 1).Value = "Timor Oriental"
$ws.Cells.Item(185, 2).Value = 24
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 20
$ws.Cells.Item(185, 5).Value = 4
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

# Row 186: Botsuana
$ws.Cells.Item(186, 1).Value = "Botsuana"
$ws.Cells.Item(186, 2).Value = 23
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 8
$ws.Cells.Item(186, 5).Value = 14
$ws.Cells.Item(186, 6).Value = 0
$ws.Cells.Item(186, 7).Value = 0
$ws.Cells.Item(186, 8).Value = 1

# Row 198: Curazao
$ws.Cells.Item(198, 1).Value = "Curazao"
$ws.Cells.Item(198, 2).Value = 16
$ws.Cells.Item(198, 3).Value = 0
$ws.Cells.Item(198, 4).Value = 13
$ws.Cells.Item(198, 5).Value = 2
$ws.Cells.Item(198, 6).Value = 0
$ws.Cells.Item(198, 7).Value = 0
$ws.Cells.Item(198, 8).Value = 1

# Row 199: Dominica
$ws.Cells.Item(199, 1).Value = "Dominica"
$ws.Cells.Item(199, 2).Value = 16
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 14
$ws.Cells.Item(199, 5).Value = 2
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# Row 205: Montserrat
$ws.Cells.Item(205, 1).Value = "Montserrat"
$ws.Cells.Item(205, 2).Value = 11
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 7
$ws.Cells.Item(205, 5).Value = 3
$ws.Cells.Item(205, 6).Value = 1
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 1

# Row 206: Seychelles
$ws.Cells.Item(206, 1).Value = "Seychelles"
$ws.Cells.Item(206, 2).Value = 11
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 8
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

# Row 211: Islas Virgenes Britanicas
$ws.Cells.Item(211, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(211, 2).Value = 7
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 3
$ws.Cells.Item(211, 5).Value = 3
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 1

# Row 212: Butan
$ws.Cells.Item(212, 1).Value = "Butan"
$ws.Cells.Item(212, 2).Value = 7
$ws.Cells.Item(212, 3).Value = 0
$ws.Cells.Item(212, 4).Value = 5
$ws.Cells.Item(212, 5).Value = 2
$ws.Cells.Item(212, 6).Value = 0
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 0

# Row 213: Bonaire, San Eustaquio y Saba
$ws.Cells.Item(213, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(213, 2).Value = 6
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 0
$ws.Cells.Item(213, 5).Value = 6
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

# Row 216: Comoras
$ws.Cells.Item(216, 1).Value = "Comoras"
$ws.Cells.Item(216, 2).Value = 4
$ws.Cells.Item(216, 3).Value = 1
$ws.Cells.Item(216, 4).Value = 0
$ws.Cells.Item(216, 5).Value = 3
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 1
$ws.Cells.Item(216, 8).Value = 1
